$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.096101641654968
$ws.Range("B1").Value = 2.067735910415649
$ws.Range("C1").Value = 9.21839427947998
$ws.Range("D1").Value = 2.420483112335205
$ws.Range("E1").Value = 1.300038814544678
